$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Finished linear actuator assembly: fill in the last missing data point
# (Achieved Motor Torque, J8) and match the plain centered style used by
# the other data cells in that row.
$ws.Range("J8").Value = 1450
$ws.Range("J8").HorizontalAlignment = -4108

# Move the selection cursor to where the user left off.
[void]$ws.Range("L8").Select()
